# "Generate Report for Handback"
#
# The localization-status report is refreshed: the two source files
# (a148332d... and eaf3e711...) have now been handed back and are in sync
# with en-US, so:
#   - the "Status" column everywhere flips from "Ready for handoff" to
#     "Handed back: in sync with en-US"
#   - the per-language sheets (zh-cn, de-de) get their "Latest Target File"
#     / "Latest Handback File" / "Latest Handback DateTime" columns filled
#     in with the handback xlf file names + timestamp, with a hyperlink on
#     the target-file cell (same link as the source-file hyperlink in col A)
#   - the columns that now hold longer text are widened so the report is
#     readable

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

$hlColor = 15570276   # RGB(100,149,237) == the workbook's existing HyperLink font color FF6495ED

$urlA148332d = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e89bcb11f78cb4912ce5ef39800eaabe9585e374/e2e/a148332d-fcc2-4bf9-9e30-1c1f6d009854.md"
$urlEaf3e711 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e89bcb11f78cb4912ce5ef39800eaabe9585e374/e2e/eaf3e711-3d40-44c3-a4d4-9772a4a9983d.md"

$mdA148332d = "a148332d-fcc2-4bf9-9e30-1c1f6d009854.md"
$mdEaf3e711 = "eaf3e711-3d40-44c3-a4d4-9772a4a9983d.md"

# ---------------------------------------------------------------------
# Overview sheet: Status columns (E, F) for both rows
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Cells.Item(2, 5).Value = $statusText
$wsOverview.Cells.Item(2, 6).Value = $statusText
$wsOverview.Cells.Item(3, 5).Value = $statusText
$wsOverview.Cells.Item(3, 6).Value = $statusText

$wsOverview.Columns.Item(5).ColumnWidth = 29.17
$wsOverview.Columns.Item(6).ColumnWidth = 29.17

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# Status column
$wsZh.Cells.Item(2, 3).Value = $statusText
$wsZh.Cells.Item(3, 3).Value = $statusText

# Row 2 (a148332d...)
$wsZh.Hyperlinks.Add($wsZh.Cells.Item(2, 9), $urlA148332d, "", "", $mdA148332d)
$wsZh.Cells.Item(2, 9).Font.Underline = 2
$wsZh.Cells.Item(2, 9).Font.Color = $hlColor
$wsZh.Cells.Item(2, 10).Value = "a148332d-fcc2-4bf9-9e30-1c1f6d009854.2f9b05af06dea172417f4dbdce4b686a95a478d2.zh-cn.xlf"
$wsZh.Cells.Item(2, 11).Value = "2016-08-21 06:53:30"

# Row 3 (eaf3e711...)
$wsZh.Hyperlinks.Add($wsZh.Cells.Item(3, 9), $urlEaf3e711, "", "", $mdEaf3e711)
$wsZh.Cells.Item(3, 9).Font.Underline = 2
$wsZh.Cells.Item(3, 9).Font.Color = $hlColor
$wsZh.Cells.Item(3, 10).Value = "eaf3e711-3d40-44c3-a4d4-9772a4a9983d.e8475997bc52974f2d0b0894fc16acf3ee9196d5.zh-cn.xlf"
$wsZh.Cells.Item(3, 11).Value = "2016-08-21 06:53:30"

$wsZh.Columns.Item(3).ColumnWidth = 29.17
$wsZh.Columns.Item(9).ColumnWidth = 39.17
$wsZh.Columns.Item(10).ColumnWidth = 39.17

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

# Status column
$wsDe.Cells.Item(2, 3).Value = $statusText
$wsDe.Cells.Item(3, 3).Value = $statusText

# Row 2 (a148332d...)
$wsDe.Hyperlinks.Add($wsDe.Cells.Item(2, 9), $urlA148332d, "", "", $mdA148332d)
$wsDe.Cells.Item(2, 9).Font.Underline = 2
$wsDe.Cells.Item(2, 9).Font.Color = $hlColor
$wsDe.Cells.Item(2, 10).Value = "a148332d-fcc2-4bf9-9e30-1c1f6d009854.2f9b05af06dea172417f4dbdce4b686a95a478d2.de-de.xlf"
$wsDe.Cells.Item(2, 11).Value = "2016-08-21 06:53:37"

# Row 3 (eaf3e711...)
$wsDe.Hyperlinks.Add($wsDe.Cells.Item(3, 9), $urlEaf3e711, "", "", $mdEaf3e711)
$wsDe.Cells.Item(3, 9).Font.Underline = 2
$wsDe.Cells.Item(3, 9).Font.Color = $hlColor
$wsDe.Cells.Item(3, 10).Value = "eaf3e711-3d40-44c3-a4d4-9772a4a9983d.e8475997bc52974f2d0b0894fc16acf3ee9196d5.de-de.xlf"
$wsDe.Cells.Item(3, 11).Value = "2016-08-21 06:53:37"

$wsDe.Columns.Item(3).ColumnWidth = 29.17
$wsDe.Columns.Item(9).ColumnWidth = 39.17
$wsDe.Columns.Item(10).ColumnWidth = 39.17
